# Weekly refresh of fruit/vegetable price data.
# Each data row (2-35) gets new values for: Fecha (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M),
# and Precio $/Kg (P). Columns A,B,C,E,F,G,H,I,N,O,Q,R are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: RowNumber, D(Fecha), J(Volumen), K(Precio minimo), L(Precio maximo), M(Precio promedio ponderado), P(Precio $/Kg)
$rows = @(
    @(2, 44455, 35, 22000, 22000, 22000, 1467),
    @(3, 44336, 65, 14000, 15000, 14462, 964),
    @(4, 44308, 40, 16000, 16000, 16000, 1067),
    @(5, 44411, 50, 22000, 22000, 22000, 1467),
    @(6, 44320, 40, 15000, 15000, 15000, 1000),
    @(7, 44321, 38, 15000, 15000, 15000, 1000),
    @(8, 44315, 65, 14000, 15000, 14538, 969),
    @(9, 44448, 85, 21000, 22000, 21529, 1435),
    @(10, 44344, 40, 20000, 20000, 20000, 1333),
    @(11, 44313, 40, 14000, 14000, 14000, 933),
    @(12, 44327, 35, 15000, 15000, 15000, 1000),
    @(13, 44322, 70, 14000, 15000, 14500, 967),
    @(14, 44314, 45, 15000, 15000, 15000, 1000),
    @(15, 44328, 38, 15000, 15000, 15000, 1000),
    @(16, 44377, 80, 18000, 19000, 18500, 1233),
    @(17, 44399, 38, 22000, 22000, 22000, 1467),
    @(18, 44316, 45, 14000, 15000, 14444, 963),
    @(19, 44397, 73, 21000, 22000, 21521, 1435),
    @(20, 44319, 50, 15000, 15000, 15000, 1000),
    @(21, 44329, 35, 15000, 15000, 15000, 1000),
    @(22, 44323, 40, 15000, 15000, 15000, 1000),
    @(23, 44334, 50, 14000, 14000, 14000, 933),
    @(24, 44333, 35, 15000, 15000, 15000, 1000),
    @(25, 44340, 47, 14000, 14000, 14000, 933),
    @(26, 44312, 80, 13000, 14000, 13562, 904),
    @(27, 44341, 40, 15000, 15000, 15000, 1000),
    @(28, 44370, 50, 18000, 18000, 18000, 1200),
    @(29, 44452, 73, 22000, 23000, 22479, 1499),
    @(30, 44326, 45, 15000, 15000, 15000, 1000),
    @(31, 44406, 50, 22000, 22000, 22000, 1467),
    @(32, 44309, 50, 15000, 15000, 15000, 1000),
    @(33, 44330, 30, 15000, 15000, 15000, 1000),
    @(34, 44343, 40, 15000, 15000, 15000, 1000),
    @(35, 44438, 75, 19000, 20000, 19467, 1298)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 4).Value = $r[1]   # D - Fecha
    $ws.Cells.Item($rowNum, 10).Value = $r[2]  # J - Volumen
    $ws.Cells.Item($rowNum, 11).Value = $r[3]  # K - Precio minimo
    $ws.Cells.Item($rowNum, 12).Value = $r[4]  # L - Precio maximo
    $ws.Cells.Item($rowNum, 13).Value = $r[5]  # M - Precio promedio ponderado
    $ws.Cells.Item($rowNum, 16).Value = $r[6]  # P - Precio $/Kg
}
